# Update the crypto price/volume table (columns D & E, rows 2-51) with the
# latest scraped values, matching the GitHub Actions "Updated cryptos list" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be stored as plain text so that values such
# as "22.00", "0.00000000363" or "30.131.22" are preserved exactly as scraped,
# instead of being auto-converted to numbers (which would drop trailing zeros,
# switch to scientific notation, or introduce floating point noise).
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @(
    @{ Row = 2; Price = '30.131.22'; Volume = '  -1.46%  ' }
    @{ Row = 3; Price = '2.105.01'; Volume = '  -0.30%  ' }
    @{ Row = 4; Price = '1.006'; Volume = '  -0.64%  ' }
    @{ Row = 5; Price = '348.06'; Volume = '  +3.59%  ' }
    @{ Row = 6; Price = '1.005'; Volume = '  -0.59%  ' }
    @{ Row = 7; Price = '0.5167'; Volume = '  -1.69%  ' }
    @{ Row = 8; Price = '0.4476'; Volume = '  -2.32%  ' }
    @{ Row = 9; Price = '52.47'; Volume = '  -5.09%  ' }
    @{ Row = 10; Price = '0.08964'; Volume = '  +0.05%  ' }
    @{ Row = 11; Price = '1.174'; Volume = '  +0.09%  ' }
    @{ Row = 12; Price = '25.48'; Volume = '  +4.33%  ' }
    @{ Row = 13; Price = '2.092.04'; Volume = '  -0.98%  ' }
    @{ Row = 14; Price = '6.741'; Volume = '  -1.77%  ' }
    @{ Row = 15; Price = '8.003'; Volume = '  -2.54%  ' }
    @{ Row = 16; Price = '99.61'; Volume = '  +2.65%  ' }
    @{ Row = 17; Price = '0.00001147'; Volume = '  -2.39%  ' }
    @{ Row = 18; Price = '1.006'; Volume = '  -0.55%  ' }
    @{ Row = 19; Price = '0.06683'; Volume = '  +0.02%  ' }
    @{ Row = 20; Price = '20.35'; Volume = '  +5.96%  ' }
    @{ Row = 21; Price = '1.005'; Volume = '  -0.64%  ' }
    @{ Row = 22; Price = '6.192'; Volume = '  -0.74%  ' }
    @{ Row = 23; Price = '30.229.45'; Volume = '  -1.34%  ' }
    @{ Row = 24; Price = '12.85'; Volume = '  +0.69%  ' }
    @{ Row = 25; Price = '2.355'; Volume = '  -0.39%  ' }
    @{ Row = 26; Price = '2.352.16'; Volume = '  -0.36%  ' }
    @{ Row = 27; Price = '22.00'; Volume = '  -1.47%  ' }
    @{ Row = 28; Price = '2.553'; Volume = '  +1.06%  ' }
    @{ Row = 29; Price = '163.67'; Volume = '  +0.08%  ' }
    @{ Row = 30; Price = '133.38'; Volume = '  -0.74%  ' }
    @{ Row = 31; Price = '1.185'; Volume = '  -3.44%  ' }
    @{ Row = 32; Price = '0.1066'; Volume = '  -0.43%  ' }
    @{ Row = 33; Price = '1.655'; Volume = '  +1.88%  ' }
    @{ Row = 34; Price = '6.262'; Volume = '  -0.89%  ' }
    @{ Row = 35; Price = '3.959'; Volume = '  -0.12%  ' }
    @{ Row = 36; Price = '5.967'; Volume = '  +1.20%  ' }
    @{ Row = 37; Price = '10.14'; Volume = '  -3.19%  ' }
    @{ Row = 38; Price = '0.02595'; Volume = '  +0.16%  ' }
    @{ Row = 39; Price = '0.06834'; Volume = '  +0.00%  ' }
    @{ Row = 40; Price = '0.2319'; Volume = '  +0.37%  ' }
    @{ Row = 41; Price = '12.51'; Volume = '  -0.46%  ' }
    @{ Row = 42; Price = '0.6834'; Volume = '  -0.35%  ' }
    @{ Row = 43; Price = '1.253'; Volume = '  -0.15%  ' }
    @{ Row = 44; Price = '14.36'; Volume = '  +2.62%  ' }
    @{ Row = 45; Price = '0.6411'; Volume = '  -0.51%  ' }
    @{ Row = 46; Price = '2.285'; Volume = '  -1.59%  ' }
    @{ Row = 47; Price = '0.00000000363'; Volume = '  +4.76%  ' }
    @{ Row = 48; Price = '3.666'; Volume = '  -0.41%  ' }
    @{ Row = 49; Price = '1.225'; Volume = '  -1.95%  ' }
    @{ Row = 50; Price = '83.00'; Volume = '  -0.13%  ' }
    @{ Row = 51; Price = '0.07226'; Volume = '  +0.71%  ' }
)

foreach ($u in $updates) {
    $ws.Range("D" + $u.Row).Value = $u.Price
    $ws.Range("E" + $u.Row).Value = $u.Volume
}
